$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-05-22"

# Update header cell text for the 2022 YTD column (column I)
$ws.Range("I1").Value = "2022 (through 05-22)"

# Update June value (row 6) for the 2022 column
$ws.Range("I6").Value = 80

# Update Total row (row 14) for the 2022 column
$ws.Range("I14").Value = 632
